$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "45/200 "
$ws.Range("C8").Value = 2000

$ws.Range("D8").Select()
